$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(6, 15).Value = 0.0
$ws.Cells.Item(7, 15).Value = 0.0009999275207519531
$ws.Cells.Item(10, 15).Value = 0.0
$ws.Cells.Item(11, 15).Value = 0.01605916023254395
$ws.Cells.Item(12, 15).Value = 0.02594470977783203
$ws.Cells.Item(13, 15).Value = 0.005932807922363281
$ws.Cells.Item(14, 15).Value = 0.0774693489074707
$ws.Cells.Item(17, 15).Value = 0.04897141456604004
$ws.Cells.Item(18, 15).Value = 0.0
$ws.Cells.Item(19, 15).Value = 0.002053976058959961
$ws.Cells.Item(20, 15).Value = 0.02978205680847168
$ws.Cells.Item(21, 15).Value = 0.006841659545898438
$ws.Cells.Item(23, 15).Value = 0.01694416999816895
$ws.Cells.Item(25, 15).Value = 0.005001544952392578
$ws.Cells.Item(26, 15).Value = 0.0167233943939209
$ws.Cells.Item(27, 15).Value = 0.0
$ws.Cells.Item(29, 15).Value = 0.05733132362365723
$ws.Cells.Item(30, 15).Value = 0.0213463306427002
$ws.Cells.Item(31, 15).Value = 0.02427220344543457
$ws.Cells.Item(32, 15).Value = 0.01630258560180664
$ws.Cells.Item(33, 15).Value = 0.0009999275207519531
$ws.Cells.Item(34, 15).Value = 0.01202583312988281
$ws.Cells.Item(35, 15).Value = 0.0
$ws.Cells.Item(36, 15).Value = 0.002997636795043945
$ws.Cells.Item(37, 15).Value = 0.02444839477539062
$ws.Cells.Item(38, 15).Value = 0.01499128341674805
$ws.Cells.Item(41, 15).Value = 0.01192498207092285
$ws.Cells.Item(44, 15).Value = 0.0
$ws.Cells.Item(48, 15).Value = 0.01692819595336914
$ws.Cells.Item(49, 15).Value = 0.009581327438354492
$ws.Cells.Item(50, 15).Value = 0.00203394889831543
$ws.Cells.Item(53, 15).Value = 0.01462578773498535
$ws.Cells.Item(54, 15).Value = 0.0009992122650146484
$ws.Cells.Item(57, 15).Value = 0.003154993057250977
$ws.Cells.Item(58, 15).Value = 0.0008606910705566406
$ws.Cells.Item(59, 15).Value = 0.0
$ws.Cells.Item(60, 15).Value = 0.005997180938720703
$ws.Cells.Item(61, 15).Value = 0.000682830810546875
$ws.Cells.Item(62, 15).Value = 0.002043008804321289
$ws.Cells.Item(63, 15).Value = 0.001001119613647461
$ws.Cells.Item(64, 15).Value = 0.01990151405334473
$ws.Cells.Item(65, 15).Value = 0.1624987125396729
$ws.Cells.Item(66, 15).Value = 0.03593325614929199
$ws.Cells.Item(67, 15).Value = 0.003111124038696289
$ws.Cells.Item(68, 15).Value = 0.03015756607055664
$ws.Cells.Item(69, 15).Value = 0.0
$ws.Cells.Item(70, 15).Value = 0.01012277603149414
$ws.Cells.Item(71, 15).Value = 1.693661689758301
$ws.Cells.Item(72, 15).Value = 3.157593727111816
$ws.Cells.Item(73, 15).Value = 0.04002022743225098
$ws.Cells.Item(75, 15).Value = 0.0
$ws.Cells.Item(76, 15).Value = 0.0009951591491699219
$ws.Cells.Item(78, 15).Value = 0.02655601501464844
$ws.Cells.Item(79, 15).Value = 0.01642394065856934
$ws.Cells.Item(80, 15).Value = 0.003007173538208008
$ws.Cells.Item(81, 15).Value = 0.2855050563812256
$ws.Cells.Item(82, 15).Value = 0.04378747940063477
$ws.Cells.Item(83, 15).Value = 0.1323230266571045
$ws.Cells.Item(84, 15).Value = 13.45798468589783
$ws.Cells.Item(85, 15).Value = 0.1207764148712158
$ws.Cells.Item(86, 15).Value = 0.01582813262939453
$ws.Cells.Item(88, 15).Value = 0.0159459114074707
$ws.Cells.Item(90, 15).Value = 0.1421191692352295
$ws.Cells.Item(91, 15).Value = 0.008006572723388672
$ws.Cells.Item(92, 15).Value = 0.1319520473480225
$ws.Cells.Item(94, 15).Value = 0.0
$ws.Cells.Item(95, 15).Value = 0.001003026962280273
$ws.Cells.Item(96, 15).Value = 0.07441973686218262
$ws.Cells.Item(97, 15).Value = 0.0
$ws.Cells.Item(98, 15).Value = 0.002004861831665039
$ws.Cells.Item(99, 15).Value = 0.001997709274291992
$ws.Cells.Item(100, 15).Value = 0.05442190170288086
$ws.Cells.Item(103, 15).Value = 0.0009992122650146484
$ws.Cells.Item(104, 15).Value = 0.02299046516418457
$ws.Cells.Item(105, 15).Value = 0.0
$ws.Cells.Item(106, 15).Value = 0.01633620262145996
$ws.Cells.Item(107, 15).Value = 0.03007650375366211
$ws.Cells.Item(108, 15).Value = 0.02606344223022461
$ws.Cells.Item(109, 15).Value = 0.001001358032226562
$ws.Cells.Item(110, 15).Value = 0.001001596450805664
$ws.Cells.Item(111, 15).Value = 0.002002954483032227
$ws.Cells.Item(112, 15).Value = 0.001991033554077148
$ws.Cells.Item(113, 15).Value = 0.006363391876220703
$ws.Cells.Item(114, 15).Value = 0.001998662948608398
$ws.Cells.Item(115, 15).Value = 0.003018617630004883
$ws.Cells.Item(116, 15).Value = 0.03776764869689941
$ws.Cells.Item(117, 15).Value = 0.001088857650756836
$ws.Cells.Item(118, 15).Value = 0.0005028247833251953
$ws.Cells.Item(119, 15).Value = 0.0
$ws.Cells.Item(120, 15).Value = 0.00100398063659668
$ws.Cells.Item(122, 15).Value = 0.0145106315612793
$ws.Cells.Item(123, 15).Value = 0.0
$ws.Cells.Item(126, 15).Value = 0.0
$ws.Cells.Item(127, 15).Value = 0.0005195140838623047
$ws.Cells.Item(128, 15).Value = 0.0009992122650146484
$ws.Cells.Item(129, 15).Value = 0.0009965896606445312
$ws.Cells.Item(130, 15).Value = 0.0
$ws.Cells.Item(132, 15).Value = 0.0007159709930419922
$ws.Cells.Item(133, 15).Value = 0.0
$ws.Cells.Item(136, 15).Value = 0.001001358032226562
$ws.Cells.Item(140, 15).Value = 0.001006126403808594
$ws.Cells.Item(141, 15).Value = 0.001611709594726562
$ws.Cells.Item(143, 15).Value = 0.0
$ws.Cells.Item(146, 15).Value = 0.0008375644683837891
$ws.Cells.Item(147, 15).Value = 0.0
$ws.Cells.Item(148, 15).Value = 0.00602269172668457
$ws.Cells.Item(149, 15).Value = 0.01870250701904297
$ws.Cells.Item(150, 15).Value = 0.02051353454589844
$ws.Cells.Item(151, 15).Value = 0.05555319786071777
$ws.Cells.Item(152, 15).Value = 0.0009152889251708984
$ws.Cells.Item(153, 15).Value = 0.0
$ws.Cells.Item(155, 15).Value = 0.0009999275207519531
$ws.Cells.Item(159, 15).Value = 0.04638528823852539
$ws.Cells.Item(160, 15).Value = 0.0
$ws.Cells.Item(161, 15).Value = 0.008521795272827148
$ws.Cells.Item(163, 15).Value = 0.001996040344238281
$ws.Cells.Item(164, 15).Value = 0.0009965896606445312
$ws.Cells.Item(166, 15).Value = 0.002916574478149414
$ws.Cells.Item(167, 15).Value = 0.001841306686401367
$ws.Cells.Item(168, 15).Value = 0.0
$ws.Cells.Item(169, 15).Value = 0.009624958038330078
$ws.Cells.Item(170, 15).Value = 0.01565456390380859
$ws.Cells.Item(171, 15).Value = 0.04068517684936523
$ws.Cells.Item(172, 15).Value = 0.009683609008789062
$ws.Cells.Item(173, 15).Value = 0.0
$ws.Cells.Item(174, 15).Value = 0.0101017951965332
$ws.Cells.Item(175, 15).Value = 0.0
$ws.Cells.Item(176, 15).Value = 0.0
$ws.Cells.Item(177, 15).Value = 0.0
